$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (row 1) so the data that used to start at row 2
# becomes row 1; this mirrors selecting the row header and deleting it.
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Delete()
